$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 735, shifting all subsequent rows down by two.
$ws.Rows("735:736").Insert()

# Fill in the two newly inserted rows with the new price-report entries.

# New row 735
$ws.Range("A735").Value = 5
$ws.Range("B735").Value = "Macroferia Regional de Talca"
$ws.Range("C735").Value = "Maule"
$ws.Range("D735").Value = 45041
$ws.Range("E735").Value = 7
$ws.Range("F735").Value = 100114001
$ws.Range("G735").Value = "Papa"
$ws.Range("H735").Value = "Asterix"
$ws.Range("I735").Value = "1a (cosecha lavada)"
$ws.Range("J735").Value = 1600
$ws.Range("K735").Value = 12000
$ws.Range("L735").Value = 12000
$ws.Range("M735").Value = 12000
$ws.Range("N735").Value = "$/malla 25 kilos"
$ws.Range("O735").Value = "Región de Los Lagos"
$ws.Range("P735").Value = 480
$ws.Range("Q735").Value = 25
$ws.Range("R735").Value = "Hortaliza"

# New row 736
$ws.Range("A736").Value = 5
$ws.Range("B736").Value = "Macroferia Regional de Talca"
$ws.Range("C736").Value = "Maule"
$ws.Range("D736").Value = 45041
$ws.Range("E736").Value = 7
$ws.Range("F736").Value = 100114001
$ws.Range("G736").Value = "Papa"
$ws.Range("H736").Value = "Patagonia"
$ws.Range("I736").Value = "1a (cosecha)"
$ws.Range("J736").Value = 1600
$ws.Range("K736").Value = 10000
$ws.Range("L736").Value = 10000
$ws.Range("M736").Value = 10000
$ws.Range("N736").Value = "$/saco 25 kilos"
$ws.Range("O736").Value = "Región de Los Lagos"
$ws.Range("P736").Value = 400
$ws.Range("Q736").Value = 25
$ws.Range("R736").Value = "Hortaliza"
